$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.621.49"
$ws.Range("E2").Value = "  +3.67%  "
$ws.Range("D3").Value = "1.696.42"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3954"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4019"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08778"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001318"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.608"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.80%  "
$ws.Range("D17").Value = "1.693.43"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "101.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07013"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.896"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "24.613.92"
$ws.Range("E24").Value = "  +3.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.069"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.336"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.189"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.494"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.02%  "
$ws.Range("D32").Value = "1.876.35"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.403"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08504"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.966"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2741"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02753"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09019"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7710"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7212"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.530"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.224"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.84%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.341"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08041"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.38%  "
